# Re-order the column headers on row 2 of Sheet1 so that the key / business-key
# columns (MilestoneID, ActivityBusinessKey, BusinessKey, MilestoneTypeBusinessKey,
# ProjectBusinessKey) come first, followed by the remaining columns in their
# previous (alphabetical) order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newHeaders = @(
    "MilestoneID",
    "ActivityBusinessKey",
    "BusinessKey",
    "MilestoneTypeBusinessKey",
    "ProjectBusinessKey",
    "Baseline",
    "BaselineDate",
    "BaselineString",
    "Code",
    "LongName",
    "Notes",
    "ReleaseDate",
    "ReportingDate",
    "ShortName",
    "Target",
    "TargetDate",
    "TargetString",
    "TextDescription",
    "UnitOfMeasure"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $newHeaders[$i]
}
